# Update StructureDefinition-dual-eligibility-indicator.xlsx
# Re-points the IG from the old "ibm.com" / "Alvearie Team" publisher to
# "linuxforhealth.org" / "LinuxForHealth Team", bumps the version/date, and
# clears the stale root-level FHIR constraint text that had been duplicated
# onto the "Extension" row.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: URL, Version, Date, Publisher ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/dual-eligibility-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet: clear the Constraint(s) cell for the root "Extension" row ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

# The "Extension.url" row's Fixed Value mirrors the StructureDefinition's own
# canonical URL (same shared string in the original workbook) - keep it in sync.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/dual-eligibility-indicator"
